$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D4"   = -7.81
    "D6"   = -8.072000000000001
    "D7"   = -8.029
    "D16"  = -8.451000000000001
    "D20"  = -7.934
    "D28"  = -7.834000000000001
    "D29"  = -7.56
    "D32"  = -7.276999999999999
    "D40"  = -7.992
    "D46"  = -8.022
    "D51"  = -8.241000000000001
    "D52"  = -7.904000000000001
    "D57"  = -7.920999999999999
    "D59"  = -8.058
    "D62"  = -7.840999999999999
    "D66"  = -7.419
    "D73"  = -8.050000000000001
    "D74"  = -8.040999999999999
    "D92"  = -7.068
    "D100" = -7.937
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
